$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: 资料小铺
$ws.Range("A10").Value = '资料小铺'
$ws.Range("B10").Formula = '=HYPERLINK("https://m.tb.cn/h.7fan38B?tk=WD5gfCUq3SV", "https://m.tb.cn/h.7fan38B?tk=WD5gfCUq3SV")'
$ws.Range("C10").Value = "https://m.tb.cn/h.7fan38B?tk=WD5gfCUq3SV"

# Row 11: 路兮
$ws.Range("A11").Value = '路兮'
$ws.Range("B11").Formula = '=HYPERLINK("https://m.tb.cn/h.7faHWPw?tk=WXOQfCUEeTG", "https://m.tb.cn/h.7faHWPw?tk=WXOQfCUEeTG")'
$ws.Range("C11").Value = "【闲鱼】https://m.tb.cn/h.7faHWPw?tk=WXOQfCUEeTG HU071 「路兮的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 12: tb_136514646
$ws.Range("A12").Value = 'tb_136514646'
$ws.Range("B12").Formula = '=HYPERLINK("https://m.tb.cn/h.7UeyCOY?tk=2YpTfC5ZI4E", "https://m.tb.cn/h.7UeyCOY?tk=2YpTfC5ZI4E")'
$ws.Range("C12").Value = "【闲鱼】https://m.tb.cn/h.7UeyCOY?tk=2YpTfC5ZI4E CZ193 「tb_136514646的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 13: 空谷画兰
$ws.Range("A13").Value = '空谷画兰'
$ws.Range("B13").Formula = '=HYPERLINK("https://m.tb.cn/h.74NtoYK?tk=JWjgfC5dr1q", "https://m.tb.cn/h.74NtoYK?tk=JWjgfC5dr1q")'
$ws.Range("C13").Value = "【闲鱼】https://m.tb.cn/h.74NtoYK?tk=JWjgfC5dr1q HU591 「空谷画兰的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 14: 心妍映雪
$ws.Range("A14").Value = '心妍映雪'
$ws.Range("B14").Formula = '=HYPERLINK("https://m.tb.cn/h.74NuTbH?tk=z2rNfC5Vrjp", "https://m.tb.cn/h.74NuTbH?tk=z2rNfC5Vrjp")'
$ws.Range("C14").Value = "【闲鱼】https://m.tb.cn/h.74NuTbH?tk=z2rNfC5Vrjp CZ225 「心妍映雪的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 15: BC素材铺
$ws.Range("A15").Value = 'BC素材铺'
$ws.Range("B15").Formula = '=HYPERLINK("https://m.tb.cn/h.74NxfEu?tk=eg30fC55tdj", "https://m.tb.cn/h.74NxfEu?tk=eg30fC55tdj")'
$ws.Range("C15").Value = "【闲鱼】https://m.tb.cn/h.74NxfEu?tk=eg30fC55tdj HU926 「BC素材铺的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 16: 是啊花呀
$ws.Range("A16").Value = '是啊花呀'
$ws.Range("B16").Formula = '=HYPERLINK("https://m.tb.cn/h.7UVgR1r?tk=A2d0fC5S5AG", "https://m.tb.cn/h.7UVgR1r?tk=A2d0fC5S5AG")'
$ws.Range("C16").Value = "【闲鱼】https://m.tb.cn/h.7UVgR1r?tk=A2d0fC5S5AG CZ193 「是啊花呀的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 17: 小狗资料铺
$ws.Range("A17").Value = '小狗资料铺'
$ws.Range("B17").Formula = '=HYPERLINK("https://m.tb.cn/h.7fZgbLp?tk=VucMfC5Rh4q", "https://m.tb.cn/h.7fZgbLp?tk=VucMfC5Rh4q")'
$ws.Range("C17").Value = "【闲鱼】https://m.tb.cn/h.7fZgbLp?tk=VucMfC5Rh4q CZ007 「小狗资料铺的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 18: 知二素材
$ws.Range("A18").Value = '知二素材'
$ws.Range("B18").Formula = '=HYPERLINK("https://m.tb.cn/h.74ndtH2?tk=INhpfC5lNVx", "https://m.tb.cn/h.74ndtH2?tk=INhpfC5lNVx")'
$ws.Range("C18").Value = "【闲鱼】https://m.tb.cn/h.74ndtH2?tk=INhpfC5lNVx CZ356 「知二素材的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 19: 咕咕资源库
$ws.Range("A19").Value = '咕咕资源库'
$ws.Range("B19").Formula = '=HYPERLINK("https://m.tb.cn/h.7fhgjkr?tk=mvsCfC5KhNt", "https://m.tb.cn/h.7fhgjkr?tk=mvsCfC5KhNt")'
$ws.Range("C19").Value = "【闲鱼】https://m.tb.cn/h.7fhgjkr?tk=mvsCfC5KhNt HU287 「咕咕资源库的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 20: 基隆圆圆的苹果
$ws.Range("A20").Value = '基隆圆圆的苹果'
$ws.Range("B20").Formula = '=HYPERLINK("https://m.tb.cn/h.74MYh6O?tk=9p9mfCgjX4S", "https://m.tb.cn/h.74MYh6O?tk=9p9mfCgjX4S")'
$ws.Range("C20").Value = "【闲鱼】https://m.tb.cn/h.74MYh6O?tk=9p9mfCgjX4S HU591 「基隆圆圆的苹果的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 21: 小北轻创业
$ws.Range("A21").Value = '小北轻创业'
$ws.Range("B21").Formula = '=HYPERLINK("https://m.tb.cn/h.74LZMia?tk=dXXSfC69dLZ", "https://m.tb.cn/h.74LZMia?tk=dXXSfC69dLZ")'
$ws.Range("C21").Value = "【闲鱼】https://m.tb.cn/h.74LZMia?tk=dXXSfC69dLZ HU071 「小北轻创业的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 22: 卡布奇诺
$ws.Range("A22").Value = '卡布奇诺'
$ws.Range("B22").Formula = '=HYPERLINK("https://m.tb.cn/h.74LYLg5?tk=nbazfC6OsxM", "https://m.tb.cn/h.74LYLg5?tk=nbazfC6OsxM")'
$ws.Range("C22").Value = "【闲鱼】https://m.tb.cn/h.74LYLg5?tk=nbazfC6OsxM CZ193 「卡布奇诺的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 23: 洋芋头的店
$ws.Range("A23").Value = '洋芋头的店'
$ws.Range("B23").Formula = '=HYPERLINK("https://m.tb.cn/h.74LePlT?tk=m5LkfC6oqkG", "https://m.tb.cn/h.74LePlT?tk=m5LkfC6oqkG")'
$ws.Range("C23").Value = "【闲鱼】https://m.tb.cn/h.74LePlT?tk=m5LkfC6oqkG MF278 「洋芋头的店的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 24: 怕黑的海龟
$ws.Range("A24").Value = '怕黑的海龟'
$ws.Range("B24").Formula = '=HYPERLINK("https://m.tb.cn/h.7UUqoxG?tk=VXX7fC6sGJL", "https://m.tb.cn/h.7UUqoxG?tk=VXX7fC6sGJL")'
$ws.Range("C24").Value = "【闲鱼】https://m.tb.cn/h.7UUqoxG?tk=VXX7fC6sGJL CZ007 「怕黑的海龟的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 25: 森森资源小铺
$ws.Range("A25").Value = '森森资源小铺'
$ws.Range("B25").Formula = '=HYPERLINK("https://m.tb.cn/h.7UUtWft?tk=QC8GfC6Flq1", "https://m.tb.cn/h.7UUtWft?tk=QC8GfC6Flq1")'
$ws.Range("C25").Value = "【闲鱼】https://m.tb.cn/h.7UUtWft?tk=QC8GfC6Flq1 CZ193 「森森资源小铺的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 26: 思秋云舒
$ws.Range("A26").Value = '思秋云舒'
$ws.Range("B26").Formula = '=HYPERLINK("https://m.tb.cn/h.7fYGf0W?tk=uVG5fC6Fv0L", "https://m.tb.cn/h.7fYGf0W?tk=uVG5fC6Fv0L")'
$ws.Range("C26").Value = "【闲鱼】https://m.tb.cn/h.7fYGf0W?tk=uVG5fC6Fv0L CZ009 「思秋云舒的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 27: 悠然资料铺
$ws.Range("A27").Value = '悠然资料铺'
$ws.Range("B27").Formula = '=HYPERLINK("https://m.tb.cn/h.7UUwOQE?tk=XAXxfC6xs6E", "https://m.tb.cn/h.7UUwOQE?tk=XAXxfC6xs6E")'
$ws.Range("C27").Value = "【闲鱼】https://m.tb.cn/h.7UUwOQE?tk=XAXxfC6xs6E CZ057 「悠然资料铺的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 28: Miao
$ws.Range("A28").Value = 'Miao'
$ws.Range("B28").Formula = '=HYPERLINK("https://m.tb.cn/h.7fYyTZn?tk=7yXIfC6zguE", "https://m.tb.cn/h.7fYyTZn?tk=7yXIfC6zguE")'
$ws.Range("C28").Value = "【闲鱼】https://m.tb.cn/h.7fYyTZn?tk=7yXIfC6zguE MF287 「Miao的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Row 29: 资料分享吧
$ws.Range("A29").Value = '资料分享吧'
$ws.Range("B29").Formula = '=HYPERLINK("https://m.tb.cn/h.7f1Y5ed?tk=TVijfChXEep", "https://m.tb.cn/h.7f1Y5ed?tk=TVijfChXEep")'
$ws.Range("C29").Value = "【闲鱼】https://m.tb.cn/h.7f1Y5ed?tk=TVijfChXEep MF278 「资料分享吧的闲鱼号，快来关注TA吧～」`n点击链接直接打开"

# Apply hyperlink-style formatting (matching B2's style) to all new B cells, and row heights
$ws.Range("B2").Copy()
$ws.Range("B10:B29").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("10:29").RowHeight = 25.5

